$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 17:05"

# --- Update per-country statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 1710114
$ws.Cells.Item(4, 3).Value = 3888
$ws.Cells.Item(4, 4).Value = 466076
$ws.Cells.Item(4, 5).Value = 1144117
$ws.Cells.Item(4, 7).Value = 116
$ws.Cells.Item(4, 8).Value = 99921

# Row 11: Alemania
$ws.Cells.Item(11, 2).Value = 181062
$ws.Cells.Item(11, 3).Value = 273
$ws.Cells.Item(11, 5).Value = 10614
$ws.Cells.Item(11, 7).Value = 20
$ws.Cells.Item(11, 8).Value = 8448

# Row 16: Canada
$ws.Cells.Item(16, 2).Value = 85997
$ws.Cells.Item(16, 3).Value = 286
$ws.Cells.Item(16, 4).Value = 44898
$ws.Cells.Item(16, 5).Value = 34533
$ws.Cells.Item(16, 7).Value = 21
$ws.Cells.Item(16, 8).Value = 6566

# Row 45: Republica Dominicana
$ws.Cells.Item(45, 2).Value = 15264
$ws.Cells.Item(45, 3).Value = 191
$ws.Cells.Item(45, 4).Value = 8534
$ws.Cells.Item(45, 5).Value = 6262
$ws.Cells.Item(45, 7).Value = 8
$ws.Cells.Item(45, 8).Value = 468

# Row 63: Moldavia
$ws.Cells.Item(63, 2).Value = 7305
$ws.Cells.Item(63, 3).Value = 158
$ws.Cells.Item(63, 5).Value = 3154
$ws.Cells.Item(63, 7).Value = 6
$ws.Cells.Item(63, 8).Value = 267

# Row 124: Sierra Leona
$ws.Cells.Item(124, 2).Value = 754
$ws.Cells.Item(124, 3).Value = 19
$ws.Cells.Item(124, 4).Value = 297
$ws.Cells.Item(124, 5).Value = 413
$ws.Cells.Item(124, 7).Value = 2
$ws.Cells.Item(124, 8).Value = 44

# Row 150: Liberia
$ws.Cells.Item(150, 2).Value = 266
$ws.Cells.Item(150, 3).Value = 1
$ws.Cells.Item(150, 4).Value = 144
$ws.Cells.Item(150, 5).Value = 96

# Row 155: Mozambique
$ws.Cells.Item(155, 2).Value = 213
$ws.Cells.Item(155, 3).Value = 4
$ws.Cells.Item(155, 5).Value = 141

# Rows 207/208: Groenlandia and Islas Turcas y Caicos swap order/values
# Before: row207 = Islas Turcas y Caicos (D=10, H=1); row208 = Groenlandia (D=11, H=0)
# After:  row207 = Groenlandia (D=11, H=0); row208 = Islas Turcas y Caicos (D=10, H=1)
$ws.Cells.Item(207, 1).Value = "Groenlandia"
$ws.Cells.Item(207, 4).Value = 11
$ws.Cells.Item(207, 8).Value = 0

$ws.Cells.Item(208, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(208, 4).Value = 10
$ws.Cells.Item(208, 8).Value = 1
